$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 4 (item id 5470)
$ws.Range("H4").Value = 638.25
$ws.Range("I4").Value = 763.44446
$ws.Range("K4").Value = 763.44446
$ws.Range("M4").Value = -649.44446

# ALC row 9 (item id 5487)
$ws.Range("H9").Value = 1500739.6
$ws.Range("I9").Value = 115.333336
$ws.Range("K9").Value = 115.333336
$ws.Range("M9").Value = 53.666664

# ALC row 63 (item id 10652)
$ws.Range("H63").Value = 69420
$ws.Range("J63").Value = 69420
$ws.Range("L63").Value = 69420
$ws.Range("N63").Value = -70668

# ALC row 66 (item id 10652)
$ws.Range("H66").Value = 69420
$ws.Range("J66").Value = 69420
$ws.Range("L66").Value = 208260
$ws.Range("N66").Value = -214500

# ALC row 113 (item id 27775)
$ws.Range("H113").Value = 201651.5
$ws.Range("I113").Value = 2125.75
$ws.Range("K113").Value = 2125.75
$ws.Range("M113").Value = 1128.25

# ALC row 138 (item id 44169)
$ws.Range("H138").Value = 2195.5518
$ws.Range("I138").Value = 917.64
$ws.Range("K138").Value = 2752.92
$ws.Range("M138").Value = 2387.08

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32 (item id 44147)
$ws.Range("H32").Value = 39655.965
$ws.Range("I32").Value = 42393.445
$ws.Range("K32").Value = 42393.445
$ws.Range("M32").Value = -42106.445

# ARM row 34 (item id 2753)
$ws.Range("H34").Value = 213799.8
$ws.Range("J34").Value = 233333
$ws.Range("L34").Value = 233333
$ws.Range("N34").Value = -233875

# ARM row 45 (item id 27714)
$ws.Range("H45").Value = 7836.2666
$ws.Range("I45").Value = 9965.546
$ws.Range("J45").Value = 1980.75
$ws.Range("K45").Value = 9965.546
$ws.Range("L45").Value = 1980.75
$ws.Range("M45").Value = -9588.546
$ws.Range("N45").Value = -2734.75

# ARM row 74 (item id 44000)
$ws.Range("H74").Value = 532442.9399999999
$ws.Range("I74").Value = 912409.4399999999
$ws.Range("J74").Value = 9989
$ws.Range("K74").Value = 912409.4399999999
$ws.Range("L74").Value = 9989
$ws.Range("M74").Value = -911535.4399999999
$ws.Range("N74").Value = -11737

# ARM row 77 (item id 44000)
$ws.Range("H77").Value = 532442.9399999999
$ws.Range("I77").Value = 912409.4399999999
$ws.Range("J77").Value = 9989
$ws.Range("K77").Value = 4562047.199999999
$ws.Range("L77").Value = 49945
$ws.Range("M77").Value = -4557679.199999999
$ws.Range("N77").Value = -58681

# ARM row 80 (item id 10667)
$ws.Range("H80").Value = 144999.5
$ws.Range("J80").Value = 144999.5
$ws.Range("L80").Value = 144999.5
$ws.Range("N80").Value = -146995.5

# ARM row 83 (item id 10667)
$ws.Range("H83").Value = 144999.5
$ws.Range("J83").Value = 144999.5
$ws.Range("L83").Value = 434998.5
$ws.Range("N83").Value = -444982.5

# ARM row 88 (item id 12530)
$ws.Range("H88").Value = 384.2
$ws.Range("J88").Value = 369.75
$ws.Range("L88").Value = 369.75
$ws.Range("N88").Value = -1181.75

# ARM row 91 (item id 12530)
$ws.Range("H91").Value = 384.2
$ws.Range("J91").Value = 369.75
$ws.Range("L91").Value = 369.75
$ws.Range("N91").Value = -3177.75

# ARM row 102 (item id 19945)
$ws.Range("H102").Value = 8374.0625
$ws.Range("I102").Value = 8213.286
$ws.Range("J102").Value = 9499.5
$ws.Range("K102").Value = 8213.286
$ws.Range("L102").Value = 9499.5
$ws.Range("M102").Value = -6591.286
$ws.Range("N102").Value = -12743.5

# ARM row 132 (item id 43997)
$ws.Range("H132").Value = 3920.8096
$ws.Range("I132").Value = 2729.7942
$ws.Range("J132").Value = 8982.625
$ws.Range("K132").Value = 8189.382599999999
$ws.Range("L132").Value = 26947.875
$ws.Range("M132").Value = -5659.382599999999
$ws.Range("N132").Value = -32007.875

$ws = $wb.Worksheets.Item("BSM")
# BSM row 20 (item id 14149)
$ws.Range("H20").Value = 2737.8
$ws.Range("I20").Value = 1999.3334
$ws.Range("J20").Value = 3054.2856
$ws.Range("K20").Value = 1999.3334
$ws.Range("L20").Value = 3054.2856
$ws.Range("M20").Value = -1752.3334
$ws.Range("N20").Value = -3548.2856

# BSM row 86 (item id 12526)
$ws.Range("H86").Value = 502099.56
$ws.Range("I86").Value = 2779.2
$ws.Range("K86").Value = 2779.2
$ws.Range("M86").Value = -1656.2

# BSM row 89 (item id 12526)
$ws.Range("H89").Value = 502099.56
$ws.Range("I89").Value = 2779.2
$ws.Range("K89").Value = 13896
$ws.Range("M89").Value = -8280

# BSM row 139 (item id 43261)
$ws.Range("H139").Value = 83399.8
$ws.Range("I139").Value = 29999
$ws.Range("J139").Value = 96750
$ws.Range("K139").Value = 29999
$ws.Range("L139").Value = 96750
$ws.Range("M139").Value = -24859
$ws.Range("N139").Value = -107030

$ws = $wb.Worksheets.Item("CRP")
# CRP row 16 (item id 27691)
$ws.Range("H16").Value = 4507.6665
$ws.Range("I16").Value = 4507.6665
$ws.Range("K16").Value = 4507.6665
$ws.Range("M16").Value = -4220.6665

# CRP row 22 (item id 5367)
$ws.Range("H22").Value = 1682.4706
$ws.Range("J22").Value = 3217
$ws.Range("L22").Value = 3217
$ws.Range("N22").Value = -3917

# CRP row 62 (item id 12580)
$ws.Range("H62").Value = 10499.429
$ws.Range("I62").Value = 4750
$ws.Range("J62").Value = 12799.2
$ws.Range("K62").Value = 4750
$ws.Range("L62").Value = 12799.2
$ws.Range("M62").Value = -4126
$ws.Range("N62").Value = -14047.2

# CRP row 65 (item id 12580)
$ws.Range("H65").Value = 10499.429
$ws.Range("I65").Value = 4750
$ws.Range("J65").Value = 12799.2
$ws.Range("K65").Value = 23750
$ws.Range("L65").Value = 63996
$ws.Range("M65").Value = -20630
$ws.Range("N65").Value = -70236

# CRP row 106 (item id 18661)
$ws.Range("H106").Value = 54999.5
$ws.Range("J106").Value = 79999
$ws.Range("L106").Value = 79999
$ws.Range("N106").Value = -82523

# CRP row 107 (item id 27689)
$ws.Range("H107").Value = 1069.4445
$ws.Range("I107").Value = 1042.0834
$ws.Range("K107").Value = 1042.0834
$ws.Range("M107").Value = 877.9166

# CRP row 113 (item id 27691)
$ws.Range("H113").Value = 4507.6665
$ws.Range("I113").Value = 4507.6665
$ws.Range("K113").Value = 4507.6665
$ws.Range("M113").Value = -2337.6665

$ws = $wb.Worksheets.Item("CUL")
# CUL row 4 (item id 4650)
$ws.Range("H4").Value = 34120096
$ws.Range("I4").Value = 50319440
$ws.Range("K4").Value = 150958320
$ws.Range("M4").Value = -150958208

# CUL row 51 (item id 4646)
$ws.Range("H51").Value = 4362.5
$ws.Range("I51").Value = 400
$ws.Range("K51").Value = 1200
$ws.Range("M51").Value = -740

# CUL row 52 (item id 31902)
$ws.Range("H52").Value = 600
$ws.Range("J52").Value = 600
$ws.Range("L52").Value = 1800
$ws.Range("N52").Value = -2332

# CUL row 80 (item id 12890)
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 6000
$ws.Range("L80").Value = 12000
$ws.Range("M80").Value = -5064
$ws.Range("N80").Value = -13872

# CUL row 83 (item id 12890)
$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 18000
$ws.Range("L83").Value = 36000
$ws.Range("M83").Value = -13320
$ws.Range("N83").Value = -45360

# CUL row 115 (item id 27861)
$ws.Range("H115").Value = 1502.5
$ws.Range("J115").Value = 20
$ws.Range("L115").Value = 60
$ws.Range("N115").Value = -2410

# CUL row 134 (item id 44074)
$ws.Range("H134").Value = 593.5
$ws.Range("I134").Value = 593.5
$ws.Range("K134").Value = 1780.5
$ws.Range("M134").Value = 3289.5

$ws = $wb.Worksheets.Item("GSM")
# GSM row 70 (item id 14146)
$ws.Range("H70").Value = 14236.094
$ws.Range("I70").Value = 12174.956
$ws.Range("J70").Value = 19503.445
$ws.Range("K70").Value = 12174.956
$ws.Range("L70").Value = 19503.445
$ws.Range("M70").Value = -11904.956
$ws.Range("N70").Value = -20043.445

# GSM row 73 (item id 14146)
$ws.Range("H73").Value = 14236.094
$ws.Range("I73").Value = 12174.956
$ws.Range("J73").Value = 19503.445
$ws.Range("K73").Value = 12174.956
$ws.Range("L73").Value = 19503.445
$ws.Range("M73").Value = -11238.956
$ws.Range("N73").Value = -21375.445

# GSM row 80 (item id 12521)
$ws.Range("H80").Value = 4331.222
$ws.Range("I80").Value = 3497
$ws.Range("K80").Value = 3497
$ws.Range("M80").Value = -2499

# GSM row 83 (item id 12521)
$ws.Range("H83").Value = 4331.222
$ws.Range("I83").Value = 3497
$ws.Range("K83").Value = 17485
$ws.Range("M83").Value = -12493

# GSM row 102 (item id 36169)
$ws.Range("H102").Value = 1662.238
$ws.Range("I102").Value = 1806.7059
$ws.Range("K102").Value = 1806.7059
$ws.Range("M102").Value = -184.7058999999999

$ws = $wb.Worksheets.Item("LTW")
# LTW row 22 (item id 5277)
$ws.Range("H22").Value = 3306.9092
$ws.Range("I22").Value = 2069.2104
$ws.Range("J22").Value = 4986.643
$ws.Range("K22").Value = 2069.2104
$ws.Range("L22").Value = 4986.643
$ws.Range("M22").Value = -1774.2104
$ws.Range("N22").Value = -5576.643

# LTW row 27 (item id 5277)
$ws.Range("H27").Value = 3306.9092
$ws.Range("I27").Value = 2069.2104
$ws.Range("J27").Value = 4986.643
$ws.Range("K27").Value = 2069.2104
$ws.Range("L27").Value = 4986.643
$ws.Range("M27").Value = -1962.2104
$ws.Range("N27").Value = -5200.643

# LTW row 46 (item id 5282)
$ws.Range("H46").Value = 4295.729
$ws.Range("J46").Value = 4824.8047
$ws.Range("L46").Value = 4824.8047
$ws.Range("N46").Value = -5200.8047

$ws = $wb.Worksheets.Item("WVR")
# WVR row 107 (item id 27746)
$ws.Range("H107").Value = 2094.6553
$ws.Range("I107").Value = 1423.3182
$ws.Range("J107").Value = 4204.5713
$ws.Range("K107").Value = 4269.9546
$ws.Range("L107").Value = 12613.7139
$ws.Range("M107").Value = -2349.9546
$ws.Range("N107").Value = -16453.7139

# WVR row 113 (item id 27752)
$ws.Range("H113").Value = 272.2069
$ws.Range("I113").Value = 289.375
$ws.Range("J113").Value = 189.8
$ws.Range("K113").Value = 868.125
$ws.Range("L113").Value = 569.4000000000001
$ws.Range("M113").Value = 1301.875
$ws.Range("N113").Value = -4909.4

# WVR row 126 (item id 36210)
$ws.Range("H126").Value = 4536.625
$ws.Range("I126").Value = 3054.0908
$ws.Range("J126").Value = 7798.2
$ws.Range("K126").Value = 9162.2724
$ws.Range("L126").Value = 23394.6
$ws.Range("M126").Value = -6692.2724
$ws.Range("N126").Value = -28334.6

# WVR row 136 (item id 44031)
$ws.Range("H136").Value = 6597.375
$ws.Range("I136").Value = 3570.75
$ws.Range("K136").Value = 10712.25
$ws.Range("M136").Value = -8162.25
